$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (rows 2-28) from 45170 to 45174, preserving existing cell formatting.
$ws.Range("C2:C28").Value = 45174
